# Day-ahead simple scheduling update for the Road data sheet.
# The road_grade (column C) values are increased by the road_length
# (column D) value for every data row (rows 2-50).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 50; $r++) {
    $grade  = $ws.Cells.Item($r, 3).Value2
    $length = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($r, 3).Value2 = $grade + $length
}

# Move the active selection to reflect where the author finished working.
$ws.Range("J18").Select() | Out-Null
